$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("K2").Value = -1.226851851851833
$ws.Range("R2").Value = 1.626775542720574
$ws.Range("S2").Value = 1.756382654173023

# Row 3
$ws.Range("K3").Value = -1.226851851851833
$ws.Range("R3").Value = 1.203236793039155
$ws.Range("S3").Value = 1.257328254301852

# Row 4
$ws.Range("K4").Value = 18.89814814814816
$ws.Range("R4").Value = 1.967443877059447
$ws.Range("S4").Value = 2.16691042047532

# Row 5
$ws.Range("K5").Value = 18.89814814814816
$ws.Range("R5").Value = 1.349021684597804
$ws.Range("S5").Value = 1.419425825968325

# Row 6
$ws.Range("K6").Value = 18.89814814814816
